# Commit: "Updated fund information on submission"
#
# The first group of samples (Vial Number 1-6) had their Condition
# changed from "FUCRW" to "RFP", and the funding/recharge information
# in columns M:P was updated to a new fund manager. The old
# Fund.Manager.Address value was removed entirely (no replacement).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Condition column (D) for the first six data rows changed from FUCRW to RFP
$ws.Range("D2:D7").Value = "RFP"

# Recharge / fund manager information (row 2, columns M:P)
$ws.Range("M2").Value = "N5WE"
$ws.Range("N2").Value = "Kylin Sakamoto"
$ws.Range("O2").Value = "kylin@lifesci.ucla.edu"
$ws.Range("P2").Value = ""

# Active cell selection moved to F10
$ws.Range("F10").Select() | Out-Null
